$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values (rows 2-11) from 3 to 4
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = 4
}

# Update the active selection to G7
$ws.Range("G7").Select()
